$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 216, pushing the existing data (old rows
# 216-336) down to 218-338. This is the weekly "new week prepended" update:
# the newest week's two quality-grade rows (Pintón / Primera Pintón) are
# added at the top of the date-ordered block, and the previous rows shift
# down by one pair.
$ws.Rows.Item(216).EntireRow.Insert()
$ws.Rows.Item(216).EntireRow.Insert()

# New row 216: Pintón, week of 2021-09-13 (serial 44452)
$ws.Range("A216").Value = 7
$ws.Range("B216").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C216").Value = "Ñuble"
$ws.Range("D216").Value = 44452
$ws.Range("E216").Value = 16
$ws.Range("F216").Value = "Fruta"
$ws.Range("G216").Value = 100108
$ws.Range("H216").Value = "Tropicales y subtropicales"
$ws.Range("I216").Value = 100108006
$ws.Range("J216").Value = "Plátano"
$ws.Range("K216").Value = "Sin especificar"
$ws.Range("L216").Value = "Pintón"
$ws.Range("M216").Value = 80
$ws.Range("N216").Value = 20000
$ws.Range("O216").Value = 20000
$ws.Range("P216").Value = 20000
$ws.Range("Q216").Value = "`$/caja 20 kilos"
$ws.Range("R216").Value = "Ecuador"
$ws.Range("S216").Value = 1000
$ws.Range("T216").Value = 20

# New row 217: Primera Pintón, same week
$ws.Range("A217").Value = 7
$ws.Range("B217").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C217").Value = "Ñuble"
$ws.Range("D217").Value = 44452
$ws.Range("E217").Value = 16
$ws.Range("F217").Value = "Fruta"
$ws.Range("G217").Value = 100108
$ws.Range("H217").Value = "Tropicales y subtropicales"
$ws.Range("I217").Value = 100108006
$ws.Range("J217").Value = "Plátano"
$ws.Range("K217").Value = "Sin especificar"
$ws.Range("L217").Value = "Primera Pintón"
$ws.Range("M217").Value = 240
$ws.Range("N217").Value = 21000
$ws.Range("O217").Value = 22000
$ws.Range("P217").Value = 21500
$ws.Range("Q217").Value = "`$/caja 20 kilos"
$ws.Range("R217").Value = "Ecuador"
$ws.Range("S217").Value = 1075
$ws.Range("T217").Value = 20
